$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Winners")

# New rows 24-30, continuing the "Order" sequence (23-29) and the
# existing pattern of winner data.
$rows = @(
    @(23, "m0 bublê", "m0 bublê", "0xdAFf0e93f8614D42Ae3Efa26657587Cd4d4bBc21"),
    @(24, "T1m", "T1m", "walletaddresshere"),
    @(25, "T1m", "T1m", "walletaddresshere"),
    @(26, "T1m", "T1m", "walletaddresshere"),
    @(27, "T1m", "T1m", "walletaddresshere"),
    @(28, "T1m", "T1m", "walletaddresshere"),
    @(29, "Jon3", "Jon3", "test")
)

$r = 24
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
